$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "even_MAG-GUT449.fa" (row 2), shifting the remaining
# rows up. This also updates the used range/dimension automatically.
$ws.Rows.Item(2).Delete()
